# Auto-generated Excel COM-interop edit script
# Applies updated cryptocurrency price/volume figures (and the Frax/EnergySwap
# row swap) to Sheet1, matching the target OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.478.56'
$ws.Range("E2").Value = '  +0.39%  '
$ws.Range("D3").Value = '1.573.73'
$ws.Range("E3").Value = '  +0.17%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("E5").Value = '  -0.08%  '
$ws.Range("D6").Value = '''292.32'
$ws.Range("E6").Value = '  +0.28%  '
$ws.Range("D7").Value = '''0.3722'
$ws.Range("E7").Value = '  -1.10%  '
$ws.Range("D8").Value = '''49.86'
$ws.Range("E8").Value = '  +0.23%  '
$ws.Range("D9").Value = '''0.3404'
$ws.Range("E9").Value = '  -0.40%  '
$ws.Range("D10").Value = '''1.150'
$ws.Range("E10").Value = '  +0.35%  '
$ws.Range("D11").Value = '''0.07556'
$ws.Range("E11").Value = '  -0.93%  '
$ws.Range("D14").Value = '''6.047'
$ws.Range("E14").Value = '  +0.67%  '
$ws.Range("D15").Value = '''6.971'
$ws.Range("E15").Value = '  +0.55%  '
$ws.Range("D16").Value = '1.572.09'
$ws.Range("E16").Value = '  +0.30%  '
$ws.Range("D17").Value = '''0.00001125'
$ws.Range("E17").Value = '  -0.71%  '
$ws.Range("D18").Value = '''90.92'
$ws.Range("E18").Value = '  +0.72%  '
$ws.Range("D19").Value = '''0.06755'
$ws.Range("E19").Value = '  +0.31%  '
$ws.Range("D21").Value = '''6.315'
$ws.Range("E21").Value = '  +1.78%  '
$ws.Range("D22").Value = '''16.38'
$ws.Range("E22").Value = '  -1.99%  '
$ws.Range("E23").Value = '  +1.25%  '
$ws.Range("D24").Value = '22.471.73'
$ws.Range("E24").Value = '  +0.39%  '
$ws.Range("D25").Value = '''2.373'
$ws.Range("E25").Value = '  -0.91%  '
$ws.Range("D26").Value = '''2.624'
$ws.Range("E26").Value = '  -1.16%  '
$ws.Range("E27").Value = '  -0.50%  '
$ws.Range("D28").Value = '''149.12'
$ws.Range("E28").Value = '  +1.36%  '
$ws.Range("D29").Value = '''5.059'
$ws.Range("E29").Value = '  +0.78%  '
$ws.Range("D30").Value = '''125.55'
$ws.Range("E30").Value = '  -0.95%  '
$ws.Range("D31").Value = '1.746.62'
$ws.Range("E31").Value = '  +0.20%  '
$ws.Range("D32").Value = '''1.083'
$ws.Range("E32").Value = '  +9.99%  '
$ws.Range("D33").Value = '''6.224'
$ws.Range("E33").Value = '  +1.67%  '
$ws.Range("D34").Value = '''2.013'
$ws.Range("E34").Value = '  +0.19%  '
$ws.Range("D35").Value = '''9.846'
$ws.Range("E35").Value = '  -2.75%  '
$ws.Range("D36").Value = '''0.08374'
$ws.Range("E36").Value = '  -1.22%  '
$ws.Range("D37").Value = '''0.02480'
$ws.Range("E37").Value = '  -2.23%  '
$ws.Range("D38").Value = '''0.2303'
$ws.Range("E38").Value = '  -0.51%  '
$ws.Range("D39").Value = '''1.338'
$ws.Range("E39").Value = '  -1.52%  '
$ws.Range("D40").Value = '''0.06533'
$ws.Range("E40").Value = '  +0.18%  '
$ws.Range("D41").Value = '''5.447'
$ws.Range("E41").Value = '  +0.44%  '
$ws.Range("D42").Value = '''11.35'
$ws.Range("E42").Value = '  -0.38%  '
$ws.Range("D43").Value = '''0.6240'
$ws.Range("E43").Value = '  -1.90%  '

# Rows 44 and 45 swap contents: Frax <-> EnergySwap, with refreshed data
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").Value = '''14.11'
$ws.Range("E44").Value = '  +0.79%  '
$ws.Range("B45").Value = 'Frax'
$ws.Range("C45").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D45").Value = '''1.001'
$ws.Range("E45").Value = '  -0.09%  '

$ws.Range("D46").Value = '''3.818'
$ws.Range("E46").Value = '  +0.61%  '
$ws.Range("D47").Value = '''0.5852'
$ws.Range("E47").Value = '  -1.81%  '
$ws.Range("D48").Value = '''129.98'
$ws.Range("E48").Value = '  +4.41%  '
$ws.Range("D49").Value = '''2.075'
$ws.Range("E49").Value = '  -0.84%  '
$ws.Range("D50").Value = '''1.211'
$ws.Range("E50").Value = '  -5.53%  '
$ws.Range("D51").Value = '''0.07336'
$ws.Range("E51").Value = '  +0.20%  '
